$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptos list refresh (prices + 1h volume %) per upstream diff.
# D-column cells that are plain decimal numbers must be forced to Text
# format first, otherwise Excel auto-converts the assigned string into a
# numeric value (e.g. "6.00" -> 6) and the literal text formatting is lost.

$ws.Range("D2").Value = "39.655.20"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "2.170.06"
$ws.Range("E3").Value = "  +0.58%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.22"
$ws.Range("E5").Value = "  -0.73%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.623"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "63.14"
$ws.Range("E7").Value = "  -0.60%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.391"
$ws.Range("E9").Value = "  -0.76%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0863"
$ws.Range("E10").Value = "  +0.45%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.104"
$ws.Range("E11").Value = "  +0.07%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.89"
$ws.Range("E12").Value = "  -1.48%  "
$ws.Range("D13").Value = "2.491.89"
$ws.Range("E13").Value = "  +0.58%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "21.76"
$ws.Range("E14").Value = "  -2.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.808"
$ws.Range("E15").Value = "  -1.32%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.47"
$ws.Range("E16").Value = "  -1.75%  "
$ws.Range("D17").Value = "2.170.21"
$ws.Range("E17").Value = "  +1.09%  "
$ws.Range("D18").Value = "39.625.87"
$ws.Range("E18").Value = "  +0.26%  "
$ws.Range("D19").Value = "0.0₃0938"
$ws.Range("E19").Value = "  +9.85%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "71.78"
$ws.Range("E20").Value = "  -0.88%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.00"
$ws.Range("E21").Value = "  -2.62%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "227.33"
$ws.Range("E22").Value = "  -0.62%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.31"
$ws.Range("E25").Value = "  -4.69%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "170.86"
$ws.Range("E26").Value = "  -0.88%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.48"
$ws.Range("E27").Value = "  -2.17%  "
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.45"
$ws.Range("E29").Value = "  +2.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.69"
$ws.Range("E30").Value = "  +0.26%  "
$ws.Range("E31").Value = "  +3.88%  "
$ws.Range("E32").Value = "  +0.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.51"
$ws.Range("E33").Value = "  -2.82%  "
$ws.Range("E34").Value = "  -2.79%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.00"
$ws.Range("E35").Value = "  -2.42%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0614"
$ws.Range("E36").Value = "  -1.35%  "
$ws.Range("E37").Value = "  +6.39%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.39"
$ws.Range("E38").Value = "  -0.59%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.86"
$ws.Range("E40").Value = "  +16.49%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "102.36"
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0228"
$ws.Range("E42").Value = "  -1.44%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.75"
$ws.Range("E43").Value = "  -2.91%  "
$ws.Range("D44").Value = "1.515.16"
$ws.Range("E44").Value = "  -1.28%  "
$ws.Range("E45").Value = "  +0.99%  "
$ws.Range("E46").Value = "  +1.36%  "
$ws.Range("E47").Value = "  -0.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0916"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("E49").Value = "  -1.84%  "
$ws.Range("E50").Value = "  +29.99%  "
$ws.Range("D51").Value = "2.373.05"
$ws.Range("E51").Value = "  +0.53%  "
